{"js": "// Helper to wrap a body-fragment of WordprocessingML into a full OOXML\n// package payload suitable for Range/Paragraph.insertOoxml().\nfunction pkg(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' + bodyXml + '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) \"Donde se genero:\" -> \"Donde se gener\u00f3:\" (drop the spell-check\n//    proofErr wrap around the now-correctly-spelled word).\n// ---------------------------------------------------------------------\nconst pGenero = paragraphs.items[3];\npGenero.insertOoxml(\n  pkg(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Donde se </w:t></w:r>' +\n      '<w:r><w:t>gener\u00f3</w:t></w:r>' +\n      '<w:r><w:t>:</w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\n\n// ---------------------------------------------------------------------\n// 2) \"ProductoDaoFactory: genera una instancia de ProductoDTO ...\" ->\n//    \"... de Producto.dao.mongo ...\" (split the single run into two).\n// ---------------------------------------------------------------------\nconst pFactory = paragraphs.items[9];\npFactory.insertOoxml(\n  pkg(\n    '<w:p>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>ProductoDaoFactory</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\">: genera una instancia de </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>Producto</w:t></w:r>' +\n      '<w:r><w:t>.dao.mongo</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> en el caso se seleccione como </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>db</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> mongo</w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\n\n// ---------------------------------------------------------------------\n// 3) Split the \"chat\" / \".dao.mongo.js\" paragraph: keep that paragraph\n//    with only those two runs, and add a brand-new paragraph after it\n//    describing the CHAT/SOCKET folder (carrying forward the _GoBack\n//    bookmark that used to sit at the end of the old paragraph).\n// ---------------------------------------------------------------------\nconst pChatDao = paragraphs.items[18];\npChatDao.insertOoxml(\n  pkg(\n    '<w:p>' +\n      '<w:r><w:t>chat</w:t></w:r>' +\n      '<w:r><w:t>.dao.mongo.js</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Para CHAT ---- se utiliz\u00f3 la carpeta </w:t></w:r>' +\n      '<w:r>' +\n        '<w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr>' +\n        '<w:t>SOCKET</w:t>' +\n      '</w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> donde utilizamos el archivo </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>chat.mongo</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t>//ingreso de productos funciona \\u2026..ver en http://localhost:8080/home</w:t></w:r>' +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n      '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4 & 5) Re-query paragraphs (count changed by +1 above) to reliably\n// reach \"Server\" and \"Se llama a las rutas...\" which are now shifted\n// one slot later.\n// ---------------------------------------------------------------------\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\n// 4) \"Server\" heading becomes bold + underlined.\nconst pServer = paragraphs2.items[25];\npServer.insertOoxml(\n  pkg(\n    '<w:p>' +\n      '<w:pPr><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr><w:t>Server</w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\n\n// 5) Append extra detail to the final paragraph about the \"start\" method.\nconst pStart = paragraphs2.items[26];\npStart.insertOoxml(\n  pkg(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Se llama a las rutas, a trav\u00e9s del m\u00e9todo </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>start</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> de </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>controller</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> del archivo </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>producto.route</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Applies the \"funciona ingreso de productos desde home\" edit to the\n# active document using the Word COM object model.\n\nfunction New-PackageXml {\n    param([string]$BodyXml)\n\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $BodyXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"Donde se genero:\" -> \"Donde se gener\u00f3:\" (paragraph 4, 1-indexed) -\n#    drop the spell-check proofErr wrap around the now-correct word.\n# ---------------------------------------------------------------------\n$pGenero = $d.Paragraphs.Item(4)\n$xmlGenero = New-PackageXml(\n    '<w:p><w:r><w:t xml:space=\"preserve\">Donde se </w:t></w:r><w:r><w:t>gener\u00f3</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'\n)\n$pGenero.Range.InsertXML($xmlGenero)\n\n# ---------------------------------------------------------------------\n# 2) \"ProductoDaoFactory: genera una instancia de ProductoDTO ...\" ->\n#    \"... de Producto.dao.mongo ...\" (paragraph 10, split one run in two).\n# ---------------------------------------------------------------------\n$pFactory = $d.Paragraphs.Item(10)\n$xmlFactory = New-PackageXml(\n    '<w:p>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>ProductoDaoFactory</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\">: genera una instancia de </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Producto</w:t></w:r>' +\n        '<w:r><w:t>.dao.mongo</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> en el caso se seleccione como </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>db</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> mongo</w:t></w:r>' +\n    '</w:p>'\n)\n$pFactory.Range.InsertXML($xmlFactory)\n\n# ---------------------------------------------------------------------\n# 3) Split the \"chat\" / \".dao.mongo.js\" paragraph (19): keep only those\n#    two runs there, and add a brand-new paragraph right after it that\n#    describes the CHAT/SOCKET folder - carrying the _GoBack bookmark\n#    that used to sit at the end of the old paragraph into the new one.\n# ---------------------------------------------------------------------\n$pChatDao = $d.Paragraphs.Item(19)\n$xmlChatDao = New-PackageXml(\n    '<w:p>' +\n        '<w:r><w:t>chat</w:t></w:r>' +\n        '<w:r><w:t>.dao.mongo.js</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">Para CHAT ---- se utiliz\u00f3 la carpeta </w:t></w:r>' +\n        '<w:r>' +\n            '<w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr>' +\n            '<w:t>SOCKET</w:t>' +\n        '</w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> donde utilizamos el archivo </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>chat.mongo</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:t>//ingreso de productos funciona \u2026..ver en http://localhost:8080/home</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n)\n$pChatDao.Range.InsertXML($xmlChatDao)\n\n# ---------------------------------------------------------------------\n# 4 & 5) The paragraph count grew by one above, so re-fetch \"Server\" and\n# the following paragraph by their new (shifted by +1) 1-based indices.\n# ---------------------------------------------------------------------\n\n# 4) \"Server\" heading becomes bold + underlined (was paragraph 25, now 26).\n$pServer = $d.Paragraphs.Item(26)\n$xmlServer = New-PackageXml(\n    '<w:p>' +\n        '<w:pPr><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr></w:pPr>' +\n        '<w:r><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr><w:t>Server</w:t></w:r>' +\n    '</w:p>'\n)\n$pServer.Range.InsertXML($xmlServer)\n\n# 5) Append extra detail about the \"start\" method (was paragraph 26, now 27).\n$pStart = $d.Paragraphs.Item(27)\n$xmlStart = New-PackageXml(\n    '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">Se llama a las rutas, a trav\u00e9s del m\u00e9todo </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>start</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> de </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>controller</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> del archivo </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>producto.route</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>'\n)\n$pStart.Range.InsertXML($xmlStart)\n"}
